$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 2
while ($true) {
    $idCell = $ws.Cells.Item($r, 1)
    $idVal = $idCell.Value2
    if (-not $idVal) {
        break
    }

    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value2 = 46076
    }

    $r = $r + 1
}
